# Update the "用户数" (user count) column values to reflect the new csvData.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 20
$ws.Range("B3").Value = 40
$ws.Range("B4").Value = 60

# Leave the active selection on the last updated cell, as Excel would after
# the user finishes editing the data.
$ws.Range("B4").Select()
